$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.570.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.55"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3761"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.76"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3674"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.68%  "

$ws.Range("E10").Value = "  +1.16%  "

$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.08"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.683"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("E15").Value = "  +2.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.439"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.642.32"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06923"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.585"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9980"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.578.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.095"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.418"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.60"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.361"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.386"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.826.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.847"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9822"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02852"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07409"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2555"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.224"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08906"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.389"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7150"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.26"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6579"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.358"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.043"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9982"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07998"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.212"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.31%  "
